$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: shift C1->D1->E1->C1 (rotate values)
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2: C2 becomes text label, E2 becomes numeric 1
$ws.Range("C2").Value = "c__Elusimicrobia"
$ws.Range("D2").Value = "c__Elusimicrobia"
$ws.Range("E2").Value = 1

# Row 3: C3 becomes text label, E3 becomes numeric 1
$ws.Range("C3").Value = "c__Elusimicrobia"
$ws.Range("D3").Value = "c__Elusimicrobia"
$ws.Range("E3").Value = 1
